{"js": "const replacements = [\n  [\"2024-10-22 Tuesday\", \"2024-10-23 Wednesday\"],\n  [\"20\\u00d770=\", \"24\\u00d757=\"],\n  [\"75\\u00d760=\", \"13\\u00d796=\"],\n  [\"54\\u00d799=\", \"23\\u00d738=\"],\n  [\"19\\u00d763=\", \"40\\u00d759=\"],\n  [\"65\\u00d732=\", \"43\\u00d796=\"],\n  [\"46\\u00d790=\", \"33\\u00d714=\"],\n  [\"22\\u00d733=\", \"87\\u00d796=\"],\n  [\"56\\u00d756=\", \"94\\u00d756=\"],\n  [\"27\\u00d724=\", \"89\\u00d789=\"],\n  [\"61\\u00d748=\", \"97\\u00d794=\"],\n  [\"44\\u00d763=\", \"26\\u00d753=\"],\n  [\"65\\u00d773=\", \"47\\u00d731=\"],\n  [\"55\\u00d725=\", \"55\\u00d731=\"],\n  [\"84\\u00d719=\", \"19\\u00d772=\"],\n  [\"73\\u00d761=\", \"15\\u00d777=\"],\n  [\"51\\u00d763=\", \"71\\u00d729=\"],\n  [\"71\\u00d712=\", \"72\\u00d732=\"],\n  [\"38\\u00d748=\", \"62\\u00d729=\"],\n  [\"49\\u00d795=\", \"14\\u00d790=\"],\n  [\"20\\u00d739=\", \"48\\u00d711=\"],\n  [\"32\\u00d763=\", \"74\\u00d783=\"],\n  [\"57\\u00d792=\", \"67\\u00d791=\"],\n  [\"47\\u00d726=\", \"50\\u00d749=\"],\n  [\"31\\u00d772=\", \"44\\u00d730=\"],\n  [\"29\\u00d741=\", \"95\\u00d716=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-22 Tuesday\", \"2024-10-23 Wednesday\"),\n    @(\"20\u00d770=\", \"24\u00d757=\"),\n    @(\"75\u00d760=\", \"13\u00d796=\"),\n    @(\"54\u00d799=\", \"23\u00d738=\"),\n    @(\"19\u00d763=\", \"40\u00d759=\"),\n    @(\"65\u00d732=\", \"43\u00d796=\"),\n    @(\"46\u00d790=\", \"33\u00d714=\"),\n    @(\"22\u00d733=\", \"87\u00d796=\"),\n    @(\"56\u00d756=\", \"94\u00d756=\"),\n    @(\"27\u00d724=\", \"89\u00d789=\"),\n    @(\"61\u00d748=\", \"97\u00d794=\"),\n    @(\"44\u00d763=\", \"26\u00d753=\"),\n    @(\"65\u00d773=\", \"47\u00d731=\"),\n    @(\"55\u00d725=\", \"55\u00d731=\"),\n    @(\"84\u00d719=\", \"19\u00d772=\"),\n    @(\"73\u00d761=\", \"15\u00d777=\"),\n    @(\"51\u00d763=\", \"71\u00d729=\"),\n    @(\"71\u00d712=\", \"72\u00d732=\"),\n    @(\"38\u00d748=\", \"62\u00d729=\"),\n    @(\"49\u00d795=\", \"14\u00d790=\"),\n    @(\"20\u00d739=\", \"48\u00d711=\"),\n    @(\"32\u00d763=\", \"74\u00d783=\"),\n    @(\"57\u00d792=\", \"67\u00d791=\"),\n    @(\"47\u00d726=\", \"50\u00d749=\"),\n    @(\"31\u00d772=\", \"44\u00d730=\"),\n    @(\"29\u00d741=\", \"95\u00d716=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $range = $d.Content\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)   -- Replace: 2 = wdReplaceAll\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
